$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2024-12-18 23:46:39", -0.1226930397070828, -0.001828903161391997, 0.008975747528043098),
    @("2024-12-18 23:46:40", -0.1227613595774164, -0.002091892289527995, 0.01027214166207885),
    @("2024-12-18 23:46:41", -0.1216129350903793, -0.002015566730399995, 0.009804759438178507),
    @("2024-12-18 23:46:42", -0.1218764545902377, -0.001896371311469995, 0.009244920481144097),
    @("2024-12-18 23:46:43", -0.1226051998737967, -0.001937419075695996, 0.009501506120600561),
    @("2024-12-18 23:46:44", -0.1231062122562434, -0.001992385625571996, 0.00981100190871818),
    @("2024-12-18 23:46:45", -0.1219480277877301, -0.001787450487037997, 0.008719042466500054),
    @("2024-12-18 23:46:46", -0.1228491994107026, -0.001789576265209996, 0.008793920458617726),
    @("2024-12-18 23:46:47", -0.1248857822120772, -0.001851679356091995, 0.009249936991660167),
    @("2024-12-18 23:46:48", -0.1261480731496704, -0.001939342398803996, 0.009785772271463346),
    @("2024-12-18 23:46:49", -0.1254323411747464, -0.002070583894041995, 0.01038872741713642),
    @("2024-12-18 23:46:50", -0.1250386885885382, -0.001975176945131996, 0.009878941397984793),
)

$startRow = 86
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
